# Auto-generated edit script: update cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '60.339.94'

$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.63%  '

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '2.679.03'

$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.02%  '

$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '523.76'

$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.72%  '

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '146.32'

$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.17%  '

$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.39%  '

$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.07%  '

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '2.696.91'

$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.78%  '

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '6.45'

$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.99%  '

$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.53%  '

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '0.339'

$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.22%  '

$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.76%  '

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '3.150.34'

$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.07%  '

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '60.364.04'

$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.79%  '

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '21.32'

$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.25%  '

$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.63%  '

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '2.688.26'

$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.07%  '

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '350.83'

$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.10%  '

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '4.53'

$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.40%  '

$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.46%  '

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '6.33'

$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.29%  '

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '0.997'

$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.11%  '

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '63.25'

$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.58%  '

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '0.421'

$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.90%  '

$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = '  +4.57%  '

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '0.994'

$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.20%  '

$c = $ws.Cells.Item(28, 2)
$c.NumberFormat = "@"
$c.Value = 'PEPE'

$c = $ws.Cells.Item(28, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '0.0₃0820'

$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.28%  '

$c = $ws.Cells.Item(29, 2)
$c.NumberFormat = "@"
$c.Value = 'InternetComputer(DFINITY)'

$c = $ws.Cells.Item(29, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '7.36'

$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.97%  '

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '6.84'

$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value = '  +6.18%  '

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'

$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.30%  '

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '19.19'

$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.27%  '

$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.33%  '

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '147.41'

$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.97%  '

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '4.28'

$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = '  +5.02%  '

$c = $ws.Cells.Item(36, 2)
$c.NumberFormat = "@"
$c.Value = 'SuiNetwork'

$c = $ws.Cells.Item(36, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '0.968'

$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.92%  '

$c = $ws.Cells.Item(37, 2)
$c.NumberFormat = "@"
$c.Value = 'ImmutableX'

$c = $ws.Cells.Item(37, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '1.25'

$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = '  +7.83%  '

$c = $ws.Cells.Item(38, 2)
$c.NumberFormat = "@"
$c.Value = 'Stacks'

$c = $ws.Cells.Item(38, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '1.52'

$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = '  +6.54%  '

$c = $ws.Cells.Item(39, 2)
$c.NumberFormat = "@"
$c.Value = 'Fetch.AI'

$c = $ws.Cells.Item(39, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '0.873'

$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.79%  '

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '36.93'

$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.18%  '

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '3.69'

$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.30%  '

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '285.78'

$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.35%  '

$c = $ws.Cells.Item(43, 2)
$c.NumberFormat = "@"
$c.Value = 'EnergySwap'

$c = $ws.Cells.Item(43, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '20.05'

$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.69%  '

$c = $ws.Cells.Item(44, 2)
$c.NumberFormat = "@"
$c.Value = 'Stellar'

$c = $ws.Cells.Item(44, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.0989'

$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.08%  '

$c = $ws.Cells.Item(45, 2)
$c.NumberFormat = "@"
$c.Value = 'FirstDigitalUSD'

$c = $ws.Cells.Item(45, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '0.996'

$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.48%  '

$c = $ws.Cells.Item(46, 2)
$c.NumberFormat = "@"
$c.Value = 'Mantle'

$c = $ws.Cells.Item(46, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '0.611'

$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.51%  '

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '2.143.63'

$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = '  +6.44%  '

$c = $ws.Cells.Item(48, 2)
$c.NumberFormat = "@"
$c.Value = 'RenderToken'

$c = $ws.Cells.Item(48, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '4.94'

$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = '  +3.55%  '

$c = $ws.Cells.Item(49, 2)
$c.NumberFormat = "@"
$c.Value = 'Hedera'

$c = $ws.Cells.Item(49, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '0.0539'

$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.24%  '

$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.31%  '

$c = $ws.Cells.Item(51, 2)
$c.NumberFormat = "@"
$c.Value = 'WhiteBITCoin'

$c = $ws.Cells.Item(51, 3)
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '10.44'

$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.71%  '

